$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell F1: same text style as the other header cells (B1:E1, style index 1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New "time_taken" metadata column (F2:F10) - plain values, no special style
$timestamps = @(
    "2021-10-05 10:50:15.473791",
    "2021-10-05 10:50:15.473801",
    "2021-10-05 10:50:15.473805",
    "2021-10-05 10:50:15.473807",
    "2021-10-05 10:50:15.473810",
    "2021-10-05 10:50:15.473813",
    "2021-10-05 10:50:15.473816",
    "2021-10-05 10:50:15.473818",
    "2021-10-05 10:50:15.473821"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
